$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.495.17'
$ws.Range('E2').Value = '  +2.04%  '
$ws.Range('D3').Value = '1.846.09'
$ws.Range('E3').Value = '  +1.57%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.032'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +2.64%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.70'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.028'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4365'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.66%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3760'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.74%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07373'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8707'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.35'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.41%  '
$ws.Range('D12').Value = '1.854.42'
$ws.Range('E12').Value = '  -8.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.496'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.660'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07180'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.43'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.033'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009018'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('E19').Value = '  +2.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.34'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('D21').Value = '27.499.91'
$ws.Range('E21').Value = '  +1.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.230'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.32'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.03%  '
$ws.Range('D24').Value = '2.073.98'
$ws.Range('E24').Value = '  -8.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.27'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.921'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.61'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.249'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.954'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +3.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.48'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09012'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7583'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.12%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.188'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.91%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.481'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.872'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.50%  '
$ws.Range('E36').Value = '  +1.88%  '
$ws.Range('E37').Value = '  +1.89%  '
$ws.Range('E38').Value = '  +2.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05270'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5126'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.800'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.54%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1667'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.675'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.435'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '108.64'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.47'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.700'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06397'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4622'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.849'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '39.05'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.81%  '
